# Quarterly income-statement refresh for سخوز (Siman Khuzestan):
#  - drop the oldest quarter column (was column D: "فصل دوم منتهی به 1399/06")
#    so every later quarter shifts one column to the left
#  - append the newest quarter ("فصل چهارم منتهی به 1401/12", published
#    1402-02-25) as the new rightmost column M
#  - a couple of the shifted figures also change slightly because of an
#    updated read_price algorithm

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Drop the obsolete leftmost quarter column -------------------------
# This shifts columns E:M -> D:L (values, styles, column widths and the
# shared-string table all shift/prune automatically), shrinking the used
# range from B1:M28 to B1:L28.
$ws.Range("D:D").Delete()

# --- 2) Re-create column M (the new quarter) with the same formatting as
#        the new rightmost existing column (L) ----------------------------
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)   # xlPasteFormats

# Match column M's width to the other "wide" quarter columns (stored width
# 31, which needs a ColumnWidth of 31 - 5/6 to round-trip through Excel's
# character->pixel rounding the same way the sibling columns do).
$ws.Range("M:M").ColumnWidth = 30.1666666666667

# --- 3) New quarter header (row 8) and publish-date (row 9) labels --------
$ws.Cells.Item(8, 13).Value2 = "فصل چهارم منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value2 = "1402-02-25 (2)"

# The quarter that used to read "1401-10-28 (8)" was republished with a
# restated figure set, so its publish-date label changes too.
$ws.Cells.Item(9, 9).Value2 = "1402-02-25 (10)"

# --- 4) New quarter's financial figures (column M, rows 11-27) ------------
$newQuarter = @{
    11 = 6128552
    12 = -2518181
    13 = 3610371
    14 = -819858
    15 = 0
    16 = 1272044
    17 = 4062557
    18 = -17007
    19 = 340743
    20 = 4386293
    21 = -107646
    22 = 4278647
    23 = 0
    24 = 4278647
    25 = 3056
    26 = 1400000
    27 = 3056
}
foreach ($r in $newQuarter.Keys) {
    $ws.Cells.Item([int]$r, 13).Value2 = $newQuarter[$r]
}

# --- 5) Small restatements to cells that merely shifted left (column I) ---
# (byproduct of the read_price algorithm change mentioned in the commit)
$colIFixups = @{
    11 = 3393666
    13 = 1302443
    17 = 924028
    19 = 128316
}
foreach ($r in $colIFixups.Keys) {
    $ws.Cells.Item([int]$r, 9).Value2 = $colIFixups[$r]
}
